# Updates cryptos list values (Price / Volume(1h) columns) per upstream scrape refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'47.966.56"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.49%  "
$ws.Range("D3").Value = "'2.482.73"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.66%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "'317.17"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.04%  "
$ws.Range("D6").Value = "'105.27"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.88%  "
$ws.Range("D7").Value = "'0.518"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.20%  "
$ws.Range("D8").Value = "'1.00"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("D9").Value = "'0.536"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.53%  "
$ws.Range("D10").Value = "'38.79"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -5.40%  "
$ws.Range("D11").Value = "'20.13"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.81%  "
$ws.Range("D12").Value = "'0.0800"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.11%  "
$ws.Range("E13").Value = "  +0.29%  "
$ws.Range("D14").Value = "'7.08"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.93%  "
$ws.Range("D15").Value = "'2.875.83"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Value = "'2.483.41"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.88%  "
$ws.Range("D17").Value = "'0.826"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.88%  "
$ws.Range("D18").Value = "'47.911.11"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.28%  "
$ws.Range("D19").Value = "'2.95"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +8.95%  "
$ws.Range("D20").Value = "'12.70"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.34%  "
$ws.Range("D21").Value = "'6.55"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.40%  "
$ws.Range("D22").Value = "'0.0₃0927"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.42%  "
$ws.Range("D23").Value = "'70.88"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.78%  "
$ws.Range("D24").Value = "'272.90"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.51%  "
$ws.Range("D25").Value = "'2.51"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.87%  "
$ws.Range("D27").Value = "'25.65"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.28%  "
$ws.Range("E28").Value = "  +2.02%  "
$ws.Range("D29").Value = "'9.69"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.60%  "
$ws.Range("E30").Value = "  -4.30%  "
$ws.Range("D31").Value = "'34.48"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.00%  "
$ws.Range("D32").Value = "'49.34"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.64%  "
$ws.Range("E33").Value = "  -0.07%  "
$ws.Range("D34").Value = "'18.93"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -5.36%  "
$ws.Range("D35").Value = "'5.25"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.90%  "
$ws.Range("D36").Value = "'0.0769"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.18%  "
$ws.Range("D37").Value = "'1.92"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.42%  "
$ws.Range("D38").Value = "'4.55"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.18%  "
$ws.Range("E39").Value = "  -5.23%  "
$ws.Range("D40").Value = "'122.27"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.73%  "
$ws.Range("D41").Value = "'0.111"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.78%  "
$ws.Range("E42").Value = "  +1.18%  "
$ws.Range("D43").Value = "'22.00"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.25%  "
$ws.Range("E44").Value = "  +0.26%  "
$ws.Range("D45").Value = "'1.999.47"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.19%  "
$ws.Range("D46").Value = "'3.15"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.26%  "
$ws.Range("D47").Value = "'1.89"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.57%  "
$ws.Range("E48").Value = "  -2.12%  "
$ws.Range("D49").Value = "'8.88"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.14%  "
$ws.Range("D50").Value = "'5.16"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.23%  "
$ws.Range("D51").Value = "'78.46"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.77%  "
